# feat: add 2022-Q1 data
#
# The workbook tracks quarterly fund-holding snapshots, one worksheet per
# quarter, plus a trailing "总计" (totals) roll-up sheet. This adds a new
# "2022-Q1" snapshot sheet (positioned right before "总计", matching the
# existing chronological tab order) and refreshes "总计" with the
# corresponding summary row.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Re-create "总计" so the new sheet can slot in right before it and
#    pick up the sheetId the existing tab order implies (总计 moves from
#    position 5 to 6, 2022-Q1 becomes the new position 5).
# ------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete() | Out-Null

$q3 = $wb.Worksheets.Item("2021-Q3")

$q1_2022 = $wb.Worksheets.Add($null, $q3)
$q1_2022.Name = "2022-Q1"

$total = $wb.Worksheets.Add($null, $q1_2022)
$total.Name = "总计"

foreach ($s in @($q1_2022, $total)) {
    $s.PageSetup.LeftMargin = 54
    $s.PageSetup.RightMargin = 54
    $s.PageSetup.TopMargin = 72
    $s.PageSetup.BottomMargin = 72
    $s.PageSetup.HeaderMargin = 36
    $s.PageSetup.FooterMargin = 36
    $s.Outline.SummaryRow = 1
    $s.Outline.SummaryColumn = 1
}

# ------------------------------------------------------------------
# 2. Populate "2022-Q1" with the per-fund holding breakdown. Column
#    layout / header styling matches the other quarterly sheets, so
#    copy the header row + index-column formatting from "2021-Q3".
# ------------------------------------------------------------------
$q3.Range("B1:H1").Copy($q1_2022.Range("B1:H1"))
$q3.Range("A2:A3").Copy($q1_2022.Range("A2:A3"))

$q1_2022.Range("B1").Value = "基金代码"
$q1_2022.Range("C1").Value = "基金名称"
$q1_2022.Range("D1").Value = "基金规模"
$q1_2022.Range("E1").Value = "股票总仓位"
$q1_2022.Range("F1").Value = "仓位占比"
$q1_2022.Range("G1").Value = "持有市值(亿元)"
$q1_2022.Range("H1").Value = "仓位排名"

$q1_2022.Range("A2").Value = 0
$q1_2022.Range("B2").Value = "'005313"
$q1_2022.Range("C2").Value = "万家中证1000指数增强A"
$q1_2022.Range("D2").Value = "'9.01"
$q1_2022.Range("E2").Value = "'93.72"
$q1_2022.Range("F2").Value = "'1.02"
$q1_2022.Range("G2").Value = "'0.0919"
$q1_2022.Range("H2").Value = 10

$q1_2022.Range("A3").Value = 1
$q1_2022.Range("B3").Value = "'005314"
$q1_2022.Range("C3").Value = "万家中证1000指数增强C"
$q1_2022.Range("D3").Value = "'4.95"
$q1_2022.Range("E3").Value = "'93.72"
$q1_2022.Range("F3").Value = "'1.02"
$q1_2022.Range("G3").Value = "'0.0505"
$q1_2022.Range("H3").Value = 10

# ------------------------------------------------------------------
# 3. Rebuild "总计": same 3-column (日期 / 持有数量(只) / 持有市值(亿元))
#    layout as before, with a new row on top for 2022-Q1 and every
#    other quarter shifted down by one.
# ------------------------------------------------------------------
$q3.Range("B1:D1").Copy($total.Range("B1:D1"))
$q3.Range("A2:A6").Copy($total.Range("A2:A6"))

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.14

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q3"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 1.12

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q2"
$total.Range("C4").Value = 17
$total.Range("D4").Value = 5.26

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q1"
$total.Range("C5").Value = 29
$total.Range("D5").Value = 10.4

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2020-Q4"
$total.Range("C6").Value = 17
$total.Range("D6").Value = 3.77


# Restore the original active sheet/selection (the workbook opened on the
# first tab before this edit) instead of leaving "总计" selected just
# because it was the last sheet touched above.
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Host "2022-Q1 + 总计 refreshed"
